# New quarterly data point has arrived for each region in this dataset, so every
# block's series shifts up by one row (the oldest quarter, 01/04/2019, drops off
# the top) and a new quarter (01/07/2024) is appended at the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each region block spans 21 data rows; (startRow, endRow, newQuarterText, newValue)
$blocks = @(
    @{ Start = 2;  End = 22; NewC = "01/07/2024"; NewD = 50.66514405698735 },
    @{ Start = 23; End = 43; NewC = "01/07/2024"; NewD = 43.92015762871584 },
    @{ Start = 44; End = 64; NewC = "01/07/2024"; NewD = 46.78631051752922 }
)

foreach ($block in $blocks) {
    $start = $block.Start
    $end = $block.End

    # Snapshot the current (pre-edit) C/D columns for this block before overwriting.
    $cVals = @{}
    $dVals = @{}
    for ($r = $start; $r -le $end; $r++) {
        $cVals[$r] = $ws.Cells.Item($r, 3).Value2
        $dVals[$r] = $ws.Cells.Item($r, 4).Value2
    }

    # Shift every row up by one: row r takes the old row r+1's quarter/value.
    for ($r = $start; $r -le ($end - 1); $r++) {
        $ws.Cells.Item($r, 3).NumberFormat = "@"
        $ws.Cells.Item($r, 3).Value = $cVals[$r + 1]
        $ws.Cells.Item($r, 3).Style = "Normal"

        if ($null -eq $dVals[$r + 1]) {
            $ws.Cells.Item($r, 4).ClearContents()
        } else {
            $ws.Cells.Item($r, 4).Value = $dVals[$r + 1]
        }
    }

    # Append the newly published quarter at the bottom of the block.
    $ws.Cells.Item($end, 3).NumberFormat = "@"
    $ws.Cells.Item($end, 3).Value = $block.NewC
    $ws.Cells.Item($end, 3).Style = "Normal"
    $ws.Cells.Item($end, 4).Value = $block.NewD
}
